# Remove the <w:contextualSpacing w:val="0"/> paragraph-property element
# from every paragraph in the document (body text, table cells, everywhere).
#
# Word's object model does not expose ParagraphFormat.ContextualSpacing as a
# settable COM property in this host, so the element is stripped by reading
# the whole document's WordOpenXML, deleting the <w:contextualSpacing/>
# tags, and re-importing the result into the same range with InsertXML
# (the supported way to edit a Range's underlying XML per the host's own
# guidance).

$d = $word.ActiveDocument
$r = $d.Content

$xml = $r.WordOpenXML

$pattern = '<w:contextualSpacing\b[^/]*/>'
$matches = [regex]::Matches($xml, $pattern)

if ($matches.Count -gt 0) {
    $newXml = [regex]::Replace($xml, $pattern, '')
    $r.InsertXML($newXml)
}

Write-Host "Removed contextualSpacing occurrences:" $matches.Count
